# Updated symbol list — applies the Dec 17 2022 cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'236.58"
$ws.Range("D2").Style = "Normal"
$ws.Range("D4").Value = "'5.467"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Value = "'0.05645"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'6.491"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Value = "'3.354"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Value = "'1.064"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Value = "'0.7944"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Value = "'0.1395"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Value = "'0.07342"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Value = "'0.03190"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Value = "'0.02981"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Value = "'0.09256"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = "'0.001670"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = "'3.253"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = "'0.04779"
$ws.Range("D17").Style = "Normal"
$ws.Range("E18").Value = "17OneONE"
$ws.Range("D19").Value = "'0.006215"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Value = "'0.005101"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Value = "'0.001053"
$ws.Range("D21").Style = "Normal"
$ws.Range("D27").Value = "'0.0004014"
$ws.Range("D27").Style = "Normal"
$ws.Range("D40").Value = "'0.04109"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Value = "'0.006946"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Value = "'0.003503"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "41CEJICEJIBestin24h"
$ws.Range("D44").Value = "'0.009784"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "43LocalTradersLCT"
$ws.Range("D45").Value = "'0.00005436"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Value = "'0.00000000751"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Value = "'0.6758"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Value = "'0.03709"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "47BOLOBOLOWorstin24h"
$ws.Range("D49").Value = "'0.00002102"
$ws.Range("D49").Style = "Normal"
